$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# Update the Timestamp column (AK2:AK54) on the FBS sheet.
# All of these cells shared the same text, so update them all in one go.
$wsFBS.Range("AK2:AK54").Value = "2024-10-27T10:01:23.940727"

# Numeric tweaks on the FBS sheet
$wsFBS.Range("P2").Value = 9.199999999999999
$wsFBS.Range("U2").Value = -4.3

$wsFBS.Range("Q14").Value = "W"
$wsFBS.Range("Q41").Value = "ENE"
$wsFBS.Range("Q48").Value = "SW"
$wsFBS.Range("Q50").Value = "S"

$wsFBS.Range("P54").Value = 3
$wsFBS.Range("U54").Value = -3.2

# Update on the Other sheet
$wsOther.Range("S10").Value = "ESE"
